$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Columns")
$ws2 = $wb.Worksheets.Item("Reports")

# ------------------------------------------------------------------
# 1. "Columns" sheet / Table1 — add a "Totals" column, and duplicate
#    the Employee Survey field rows as a new "ESPDF" report.
# ------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Item(1)

# Append a new "Totals" column at the end of the table.
$lo1.ListColumns.Add() | Out-Null
$lo1.HeaderRowRange.Cells.Item(1, 12).Value = "Totals"

# Duplicate the 7 existing data rows (Employee Survey) as rows 9-15.
$ws1.Range("A2:K8").Copy()
$ws1.Range("A9").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# New rows belong to report "ESPDF" instead of "Employee Survey".
$ws1.Range("A9:A15").Value = "ESPDF"

# Re-apply the header-style formatting used on column B for rows 3-8
# (bold + centered) to the equivalent new rows 10-15, and the
# bold+centered+text-format style used on B8 to the new B15.
$ws1.Range("B10:B14").Font.Bold = $true
$ws1.Range("B10:B14").HorizontalAlignment = -4108
$ws1.Range("B9").Font.Bold = $true
$ws1.Range("B9").HorizontalAlignment = -4108
$ws1.Range("B15").Font.Bold = $true
$ws1.Range("B15").HorizontalAlignment = -4108
$ws1.Range("B15").NumberFormat = "@"

# Mark "Eval Number" as a totals field for the new ESPDF report.
$ws1.Range("L15").Value = "Y"

# Grow the table to include the new column + new rows.
$lo1.Resize($ws1.Range("A1:L15"))

# ------------------------------------------------------------------
# 2. "Reports" sheet / Table2 — insert a "Heading Totals" column
#    right after "Heading Type", and rename the PDF report.
# ------------------------------------------------------------------
$lo2 = $ws2.ListObjects.Item(1)
$loName2 = $lo2.Name
$loStyle2 = $lo2.TableStyle.Name

# Insert a blank worksheet column in position E (after "Heading Type").
$ws2.Columns.Item(5).Insert()
$ws2.Cells.Item(1, 5).Value = "Heading Totals"
$ws2.Cells.Item(2, 5).Value = "N"
$ws2.Cells.Item(3, 5).Value = "Y"

# Rename the PDF report from "Employee Survey PDF" to "ESPDF".
$ws2.Cells.Item(3, 1).Value = "ESPDF"

# Rebuild the ListObject so its column collection picks up the
# inserted column in the correct position.
$lo2.Unlist()
$lo2b = $ws2.ListObjects.Add(1, $ws2.Range("A1:R3"), 0, 1)
$lo2b.Name = $loName2
$lo2b.TableStyle = $loStyle2

# ------------------------------------------------------------------
# 3. View state — active sheet/selection moves from Reports!D4 to
#    Columns!L10, with Reports left selected at E3.
# ------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("E3").Select()
$ws1.Activate()
$ws1.Range("L10").Select()
